$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "58.073.65"
$ws.Range("E2").Value = "  -2.18%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.576.84"
$ws.Range("E3").Value = "  -2.45%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.00%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'517.62"
$ws.Range("E5").Value = "  -2.15%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'139.39"
$ws.Range("E6").Value = "  -3.98%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.06%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  -1.63%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = "2.589.61"
$ws.Range("E9").Value = "  -2.53%  "

# Row 10 - Toncoin
$ws.Range("D10").Value = "'6.43"
$ws.Range("E10").Value = "  -3.72%  "

# Row 11 - Dogecoin
$ws.Range("D11").Value = "'0.0997"
$ws.Range("E11").Value = "  -4.54%  "

# Row 12 - Cardano
$ws.Range("D12").Value = "'0.327"
$ws.Range("E12").Value = "  -3.38%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  +0.41%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "3.024.45"
$ws.Range("E14").Value = "  -2.78%  "

# Row 15 - WrappedBTC
$ws.Range("D15").Value = "58.041.16"
$ws.Range("E15").Value = "  -2.16%  "

# Row 16 - Avalanche
$ws.Range("D16").Value = "'20.10"
$ws.Range("E16").Value = "  -4.26%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "2.600.02"
$ws.Range("E17").Value = "  -0.60%  "

# Row 18 - ShibaInu
$ws.Range("D18").Value = "'0.0000132"
$ws.Range("E18").Value = "  -3.99%  "

# Row 19 - BitcoinCash
$ws.Range("D19").Value = "'333.48"
$ws.Range("E19").Value = "  -2.59%  "

# Row 20 - Polkadot
$ws.Range("D20").Value = "'4.29"
$ws.Range("E20").Value = "  -3.91%  "

# Row 21 - Chainlink
$ws.Range("D21").Value = "'10.10"
$ws.Range("E21").Value = "  -5.01%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "'6.36"
$ws.Range("E22").Value = "  -0.20%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  +0.08%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "'65.91"
$ws.Range("E24").Value = "  +0.22%  "

# Row 25 - Kaspa
$ws.Range("E25").Value = "  -1.55%  "

# Row 26 - Binance-PegBSC-USD
$ws.Range("D26").Value = "'0.997"
$ws.Range("E26").Value = "  +0.00%  "

# Row 27 - Polygon
$ws.Range("D27").Value = "'0.399"
$ws.Range("E27").Value = "  -4.88%  "

# Row 28 - WrappedeETH
$ws.Range("D28").Value = "2.688.93"
$ws.Range("E28").Value = "  -2.51%  "

# Row 29 - InternetComputer(DFINITY)
$ws.Range("D29").Value = "'6.94"
$ws.Range("E29").Value = "  -4.47%  "

# Row 31 - PEPE
$ws.Range("D31").Value = "0.0₃0719"
$ws.Range("E31").Value = "  -10.42%  "

# Row 32 - Aptos
$ws.Range("D32").Value = "'5.95"
$ws.Range("E32").Value = "  -8.22%  "

# Row 33 - PancakeSwap
$ws.Range("E33").Value = "  -3.70%  "

# Row 34 - EthereumClassic
$ws.Range("D34").Value = "'18.66"
$ws.Range("E34").Value = "  -2.15%  "

# Row 35 - Monero
$ws.Range("D35").Value = "'149.16"
$ws.Range("E35").Value = "  -0.56%  "

# Row 36 - NEARProtocol
$ws.Range("D36").Value = "'3.92"
$ws.Range("E36").Value = "  -6.83%  "

# Row 37 - ImmutableX
$ws.Range("E37").Value = "  -7.70%  "

# Row 38 - OKB -> Fetch.AI (swap)
$ws.Range("B38").Value = "Fetch.AI"
$ws.Range("C38").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D38").Value = "'0.842"
$ws.Range("E38").Value = "  -2.99%  "

# Row 39 - Fetch.AI -> OKB (swap)
$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").Value = "'36.19"
$ws.Range("E39").Value = "  -1.29%  "

# Row 40 - SuiNetwork
$ws.Range("D40").Value = "'0.827"
$ws.Range("E40").Value = "  -7.60%  "

# Row 41 - Stacks
$ws.Range("E41").Value = "  -4.54%  "

# Row 42 - Filecoin
$ws.Range("D42").Value = "'3.50"
$ws.Range("E42").Value = "  -4.43%  "

# Row 43 - FirstDigitalUSD
$ws.Range("E43").Value = "  -0.03%  "

# Row 44 - Bittensor
$ws.Range("D44").Value = "'275.39"
$ws.Range("E44").Value = "  +0.87%  "

# Row 45 - WhiteBITCoin
$ws.Range("E45").Value = "  +0.27%  "

# Row 46 - Mantle
$ws.Range("E46").Value = "  -2.54%  "

# Row 47 - Stellar
$ws.Range("D47").Value = "'0.0942"
$ws.Range("E47").Value = "  -3.47%  "

# Row 48 - Hedera
$ws.Range("D48").Value = "'0.0515"
$ws.Range("E48").Value = "  -4.55%  "

# Row 49 - EnergySwap
$ws.Range("D49").Value = "'18.38"
$ws.Range("E49").Value = "  -5.36%  "

# Row 50 - Maker
$ws.Range("D50").Value = "1.968.48"
$ws.Range("E50").Value = "  -3.60%  "

# Row 51 - RenderToken
$ws.Range("D51").Value = "'4.50"
$ws.Range("E51").Value = "  -5.93%  "
